$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state for rows 2-23 (column A = empadronador, column B = total_registros)
$data = @(
    @("DAVILA CORDOVA MARIBEL", 55),
    @("FARFAN MONTOYA ROSSANA ISABEL", 54),
    @("ANGIE BELÉN RODRÍGUEZ ZAVALA", 53),
    @("SANCHEZ ULLOA CESAR AUGUSTO", 51),
    @("FIORELA KEILY GUTIERREZ CRUZ", 49),
    @("GARCIA GUTIERREZ LUIS ARTURO", 48),
    @("RAMOS RAMOS HANDY JAIR", 47),
    @("CARBAJAL RAMOS JESUS MARINA", 46),
    @("DELGADO DELGADO RONI", 46),
    @("OLIVA ALVA GOSSELYN NASSIRA", 46),
    @("GONZALES VALLE SEBASTIAN", 45),
    @("VERDE LIZARRAGA DEYSI EUFEMIA", 45),
    @("JOSSY IVANA SUÁREZ ZAVALETA", 44),
    @("BAZAN TEJADA JOSE VICENTE", 43),
    @("CASTILLO QUEZADA DIEGO ALONSO", 43),
    @("DANY DARWIN VILLACORTA SAAVEDRA", 42),
    @("GUERRA CALDERON ESTHEFANY NICOLLE", 40),
    @("JAVE CHAVEZ ANGHELO MARTIN", 39),
    @("ARANEDA LOPEZ MARCO VIERI", 38),
    @("CARDENAS CAMPOJO MARY PAULA", 36),
    @("ALVITES CAMPOS SERGIO MARTIN", 17),
    @("ROCHA SIPIRAN JHORDAN ENRIQUE", 1)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row = $row + 1
}
